$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the old row 7 ("3" row that ended the first group) - this shifts everything
# up by one row, matching the new structure.
$ws.Rows.Item(7).Delete()

# Add the new test step text in column C of the new row 9 (was row 10 before deletion).
$ws.Range("C9").Value = "Toutes les informations et données sont visible"

# Update the active selection to reflect the new cursor position.
$ws.Range("C9").Select()
